# Updates Price (D) and Volume(1h) (E) columns for the cryptos worksheet
# per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "24.844.99"
$ws.Cells.Item(2, 5).Value = "  +0.70%  "
$ws.Cells.Item(3, 4).Value = "1.658.27"
$ws.Cells.Item(3, 5).Value = "  -2.15%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.72%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "316.47"
$ws.Cells.Item(5, 5).Value = "  +1.95%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9957"
$ws.Cells.Item(6, 5).Value = "  -0.62%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3631"
$ws.Cells.Item(7, 5).Value = "  -2.27%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "47.11"
$ws.Cells.Item(8, 5).Value = "  -3.82%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.3267"
$ws.Cells.Item(9, 5).Value = "  -3.87%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.139"
$ws.Cells.Item(10, 5).Value = "  -5.08%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07061"
$ws.Cells.Item(11, 5).Value = "  -4.76%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.9977"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.043"
$ws.Cells.Item(13, 5).Value = "  -3.75%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "19.56"
$ws.Cells.Item(14, 5).Value = "  -6.02%  "
$ws.Cells.Item(15, 4).Value = "1.664.26"
$ws.Cells.Item(15, 5).Value = "  -1.66%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.617"
$ws.Cells.Item(16, 5).Value = "  -4.58%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001049"
$ws.Cells.Item(17, 5).Value = "  -5.87%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.06616"
$ws.Cells.Item(18, 5).Value = "  -1.21%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.9973"
$ws.Cells.Item(19, 5).Value = "  -0.43%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "79.22"
$ws.Cells.Item(20, 5).Value = "  -4.38%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.913"
$ws.Cells.Item(21, 5).Value = "  -5.92%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "15.75"
$ws.Cells.Item(22, 5).Value = "  -7.62%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "12.57"
$ws.Cells.Item(23, 5).Value = "  -2.21%  "
$ws.Cells.Item(24, 4).Value = "24.874.70"
$ws.Cells.Item(24, 5).Value = "  +0.72%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.435"
$ws.Cells.Item(25, 5).Value = "  -0.50%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.402"
$ws.Cells.Item(26, 5).Value = "  -12.45%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "148.78"
$ws.Cells.Item(27, 5).Value = "  +0.36%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "18.63"
$ws.Cells.Item(28, 5).Value = "  -7.43%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.226"
$ws.Cells.Item(29, 5).Value = "  +0.57%  "
$ws.Cells.Item(30, 4).Value = "1.848.54"
$ws.Cells.Item(30, 5).Value = "  -1.62%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "125.90"
$ws.Cells.Item(31, 5).Value = "  -4.06%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.076"
$ws.Cells.Item(32, 5).Value = "  -3.17%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.825"
$ws.Cells.Item(33, 5).Value = "  -12.72%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.08436"
$ws.Cells.Item(34, 5).Value = "  -2.63%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.678"
$ws.Cells.Item(35, 5).Value = "  -4.15%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "12.28"
$ws.Cells.Item(36, 5).Value = "  -9.25%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.288"
$ws.Cells.Item(37, 5).Value = "  +2.65%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "5.208"
$ws.Cells.Item(38, 5).Value = "  -5.15%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.06038"
$ws.Cells.Item(39, 5).Value = "  -8.03%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.02233"
$ws.Cells.Item(40, 5).Value = "  -6.15%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.2070"
$ws.Cells.Item(41, 5).Value = "  -5.80%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "8.222"
$ws.Cells.Item(42, 5).Value = "  -8.52%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.9960"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.5928"
$ws.Cells.Item(44, 5).Value = "  -6.64%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.809"
$ws.Cells.Item(45, 5).Value = "  +0.09%  "
$ws.Cells.Item(46, 5).Value = "  -6.64%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5644"
$ws.Cells.Item(47, 5).Value = "  -6.53%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "125.13"
$ws.Cells.Item(48, 5).Value = "  -2.23%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.953"
$ws.Cells.Item(49, 5).Value = "  -6.71%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.07021"
$ws.Cells.Item(50, 5).Value = "  -3.00%  "
